$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1657.7273
$ws.Range("I18").Value = 1715.5555
$ws.Range("J18").Value = 1397.5
$ws.Range("K18").Value = 1715.5555
$ws.Range("L18").Value = 1397.5
$ws.Range("M18").Value = -1431.5555
$ws.Range("N18").Value = -1965.5
$ws.Range("H19").Value = 1904.5714
$ws.Range("J19").Value = 1811.4286
$ws.Range("L19").Value = 1811.4286
$ws.Range("N19").Value = -2161.4286
$ws.Range("H62").Value = 6153.5
$ws.Range("I62").Value = 5713.2666
$ws.Range("K62").Value = 5713.2666
$ws.Range("M62").Value = -5089.2666
$ws.Range("H65").Value = 6153.5
$ws.Range("I65").Value = 5713.2666
$ws.Range("K65").Value = 28566.333
$ws.Range("M65").Value = -25446.333
$ws.Range("H70").Value = 41668204
$ws.Range("I70").Value = 1949.3334
$ws.Range("K70").Value = 5848.0002
$ws.Range("M70").Value = -5578.0002
$ws.Range("H73").Value = 41668204
$ws.Range("I73").Value = 1949.3334
$ws.Range("K73").Value = 5848.0002
$ws.Range("M73").Value = -4912.0002
$ws.Range("H80").Value = 75011180
$ws.Range("I80").Value = 250000430
$ws.Range("J80").Value = 16681438
$ws.Range("K80").Value = 750001290
$ws.Range("L80").Value = 50044314
$ws.Range("M80").Value = -750000292
$ws.Range("N80").Value = -50046310
$ws.Range("H83").Value = 75011180
$ws.Range("I83").Value = 250000430
$ws.Range("J83").Value = 16681438
$ws.Range("K83").Value = 2250003870
$ws.Range("L83").Value = 150132942
$ws.Range("M83").Value = -2249998878
$ws.Range("N83").Value = -150142926
$ws.Range("H88").Value = 1399.091
$ws.Range("J88").Value = 1573.1428
$ws.Range("L88").Value = 1573.1428
$ws.Range("N88").Value = -2385.1428
$ws.Range("H91").Value = 1399.091
$ws.Range("J91").Value = 1573.1428
$ws.Range("L91").Value = 1573.1428
$ws.Range("N91").Value = -4381.1428
$ws.Range("H112").Value = 3654.9607
$ws.Range("J112").Value = 3724.551
$ws.Range("L112").Value = 11173.653
$ws.Range("N112").Value = -13389.653
$ws.Range("H116").Value = 22743264
$ws.Range("I116").Value = 25016816
$ws.Range("K116").Value = 25016816
$ws.Range("M116").Value = -25013374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7966.48
$ws.Range("I32").Value = 4982.414
$ws.Range("K32").Value = 4982.414
$ws.Range("M32").Value = -4695.414
$ws.Range("H45").Value = 5968.625
$ws.Range("I45").Value = 5971.2856
$ws.Range("K45").Value = 5971.2856
$ws.Range("M45").Value = -5594.2856
$ws.Range("H61").Value = 3622.6924
$ws.Range("I61").Value = 3195.4
$ws.Range("K61").Value = 3195.4
$ws.Range("M61").Value = -2983.4
$ws.Range("H136").Value = 3622.6924
$ws.Range("I136").Value = 3195.4
$ws.Range("K136").Value = 9586.200000000001
$ws.Range("M136").Value = -7036.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 40161
$ws.Range("I20").Value = 70048
$ws.Range("K20").Value = 70048
$ws.Range("M20").Value = -69801
$ws.Range("H81").Value = 201080.33
$ws.Range("I81").Value = 19709
$ws.Range("J81").Value = 237354.6
$ws.Range("K81").Value = 19709
$ws.Range("L81").Value = 237354.6
$ws.Range("M81").Value = -18648
$ws.Range("N81").Value = -239476.6
$ws.Range("H84").Value = 201080.33
$ws.Range("I84").Value = 19709
$ws.Range("J84").Value = 237354.6
$ws.Range("K84").Value = 59127
$ws.Range("L84").Value = 712063.8
$ws.Range("M84").Value = -53823
$ws.Range("N84").Value = -722671.8
$ws.Range("H86").Value = 2451.9333
$ws.Range("I86").Value = 2484.3635
$ws.Range("J86").Value = 2362.75
$ws.Range("K86").Value = 2484.3635
$ws.Range("L86").Value = 2362.75
$ws.Range("M86").Value = -1361.3635
$ws.Range("N86").Value = -4608.75
$ws.Range("H89").Value = 2451.9333
$ws.Range("I89").Value = 2484.3635
$ws.Range("J89").Value = 2362.75
$ws.Range("K89").Value = 12421.8175
$ws.Range("L89").Value = 11813.75
$ws.Range("M89").Value = -6805.817499999999
$ws.Range("N89").Value = -23045.75
$ws.Range("H140").Value = 99891.836
$ws.Range("J140").Value = 99891.836
$ws.Range("L140").Value = 99891.836
$ws.Range("N140").Value = -110251.836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1118.2727
$ws.Range("J94").Value = 1163
$ws.Range("L94").Value = 1163
$ws.Range("N94").Value = -2065
$ws.Range("H122").Value = 33334198
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H124").Value = 55774.11
$ws.Range("J124").Value = 55774.11
$ws.Range("L124").Value = 55774.11
$ws.Range("N124").Value = -60684.11
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H140").Value = 105000
$ws.Range("I140").Value = 80000
$ws.Range("J140").Value = 130000
$ws.Range("K140").Value = 80000
$ws.Range("L140").Value = 130000
$ws.Range("M140").Value = -74820
$ws.Range("N140").Value = -140360
$ws.Range("H141").Value = 495337
$ws.Range("I141").Value = 20296
$ws.Range("J141").Value = 1029758.1
$ws.Range("K141").Value = 20296
$ws.Range("L141").Value = 1029758.1
$ws.Range("M141").Value = -15116
$ws.Range("N141").Value = -1040118.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1332.65
$ws.Range("I34").Value = 841.875
$ws.Range("J34").Value = 1659.8334
$ws.Range("K34").Value = 2525.625
$ws.Range("L34").Value = 4979.5002
$ws.Range("M34").Value = -2441.625
$ws.Range("N34").Value = -5147.5002
$ws.Range("H37").Value = 95000
$ws.Range("J37").Value = 95000
$ws.Range("L37").Value = 285000
$ws.Range("N37").Value = -285224
$ws.Range("H86").Value = 1397.5
$ws.Range("J86").Value = 2095
$ws.Range("L86").Value = 6285
$ws.Range("N86").Value = -8657
$ws.Range("H89").Value = 1397.5
$ws.Range("J89").Value = 2095
$ws.Range("L89").Value = 18855
$ws.Range("N89").Value = -30711
$ws.Range("H131").Value = 1676.1708
$ws.Range("I131").Value = 1374.75
$ws.Range("J131").Value = 1869.08
$ws.Range("K131").Value = 4124.25
$ws.Range("L131").Value = 5607.24
$ws.Range("M131").Value = 915.75
$ws.Range("N131").Value = -15687.24
$ws.Range("H138").Value = 35503080
$ws.Range("J138").Value = 53252620
$ws.Range("L138").Value = 159757860
$ws.Range("N138").Value = -159768140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4499.2393
$ws.Range("I70").Value = 4497.385
$ws.Range("J70").Value = 4499.9697
$ws.Range("K70").Value = 4497.385
$ws.Range("L70").Value = 4499.9697
$ws.Range("M70").Value = -4227.385
$ws.Range("N70").Value = -5039.9697
$ws.Range("H73").Value = 4499.2393
$ws.Range("I73").Value = 4497.385
$ws.Range("J73").Value = 4499.9697
$ws.Range("K73").Value = 4497.385
$ws.Range("L73").Value = 4499.9697
$ws.Range("M73").Value = -3561.385
$ws.Range("N73").Value = -6371.9697
$ws.Range("H80").Value = 1999.25
$ws.Range("I80").Value = 1999.25
$ws.Range("K80").Value = 1999.25
$ws.Range("M80").Value = -1001.25
$ws.Range("H83").Value = 1999.25
$ws.Range("I83").Value = 1999.25
$ws.Range("K83").Value = 9996.25
$ws.Range("M83").Value = -5004.25
$ws.Range("H132").Value = 1504.5625
$ws.Range("I132").Value = 1312.3572
$ws.Range("K132").Value = 3937.0716
$ws.Range("M132").Value = -1407.0716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15098.111
$ws.Range("I7").Value = 2193
$ws.Range("J7").Value = 21550.666
$ws.Range("K7").Value = 2193
$ws.Range("L7").Value = 21550.666
$ws.Range("M7").Value = -2081
$ws.Range("N7").Value = -21774.666
$ws.Range("H40").Value = 9808106
$ws.Range("I40").Value = 12349833
$ws.Range("J40").Value = 4305.4287
$ws.Range("K40").Value = 12349833
$ws.Range("L40").Value = 4305.4287
$ws.Range("M40").Value = -12349697
$ws.Range("N40").Value = -4577.4287
$ws.Range("H93").Value = 928.5
$ws.Range("I93").Value = 584
$ws.Range("J93").Value = 1568.2858
$ws.Range("K93").Value = 584
$ws.Range("L93").Value = 1568.2858
$ws.Range("M93").Value = 664
$ws.Range("N93").Value = -4064.2858
$ws.Range("H126").Value = 15098.111
$ws.Range("I126").Value = 2193
$ws.Range("J126").Value = 21550.666
$ws.Range("K126").Value = 6579
$ws.Range("L126").Value = 64651.99800000001
$ws.Range("M126").Value = -4109
$ws.Range("N126").Value = -69591.99800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 11762.5
$ws.Range("I40").Value = 11762.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 11762.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -11613.5
$ws.Range("N40").ClearContents()
$ws.Range("H100").Value = 2973.6155
$ws.Range("I100").Value = 3075.3333
$ws.Range("J100").Value = 2744.75
$ws.Range("K100").Value = 6150.6666
$ws.Range("L100").Value = 5489.5
$ws.Range("M100").Value = -5609.6666
$ws.Range("N100").Value = -6571.5
$ws.Range("H122").Value = 23257226
$ws.Range("I122").Value = 25642356
$ws.Range("J122").Value = 2223
$ws.Range("K122").Value = 76927068
$ws.Range("L122").Value = 6669
$ws.Range("M122").Value = -76924618
$ws.Range("N122").Value = -11569
$ws.Range("H132").Value = 23698.877
$ws.Range("I132").Value = 28696.025
$ws.Range("K132").Value = 86088.07500000001
$ws.Range("M132").Value = -83558.07500000001
